$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("pt_min"), splitting the old "eta"
# column (D) into two columns: "eta_min" (D) and "eta_max" (E).
$ws.Columns("E").Insert()

# Update header row: D1 used to be "eta", now becomes "eta_min";
# the newly inserted E1 becomes "eta_max".
$ws.Range("D1").Value = "eta_min"
$ws.Range("E1").Value = "eta_max"

# Replace the single "eta" values in D2:D5 with the eta_min/eta_max pairs.
$ws.Range("D2").Value = -2.2
$ws.Range("E2").Value = -1.2

$ws.Range("D3").Value = 1.2
$ws.Range("E3").Value = 2.2

$ws.Range("D4").Value = -2.2
$ws.Range("E4").Value = -1.2

$ws.Range("D5").Value = 1.2
$ws.Range("E5").Value = 2.2

# Restore the active selection saved in the workbook view.
$ws.Range("E18").Select()
